# Fix the bug about taxing negative income.
# On the "Intermediate" sheet, row 50 (after-tax cash flow) previously
# multiplied every year's pre-tax value (row 49) by (1 - tax rate), even
# when the pre-tax value was negative - effectively "taxing" a loss as if
# it were a credit. The fix only applies the tax multiplier when the
# pre-tax value is positive; negative (loss) years pass through unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intermediate")

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($col in $cols) {
    $ref = $col + "49"
    $formula = "=IF(" + $ref + ">0," + $ref + "*(1-'Financial param'!`$I`$10/100)," + $ref + ")"
    $ws.Range($col + "50").Formula = $formula
}

# Restore the view state recorded in the saved workbook: Intermediate sheet
# is the active sheet/tab, scrolled so B56 (the IRR result) is selected.
$ws.Activate()
$ws.Range("B56").Select()

# Financial param sheet's last remembered selection.
$wsFin = $wb.Worksheets.Item("Financial param")
$wsFin.Range("I10").Select()

# Re-activate Intermediate sheet as the final tab shown (matches
# activeTab="2" / tabSelected="1" on Intermediate in the saved file).
$ws.Activate()
